$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.753.20'
$ws.Range("E2").Value = '  +2.48%  '
$ws.Range("D3").Value = '2.381.92'
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''552.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.40%  '
$ws.Range("D6").Value = '''141.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.14%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '''0.525'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.87%  '
$ws.Range("D9").Value = '2.383.07'
$ws.Range("E9").Value = '  +1.62%  '
$ws.Range("E10").Value = '  +4.94%  '
$ws.Range("E11").Value = '  +2.30%  '
$ws.Range("E12").Value = '  +2.38%  '
$ws.Range("E13").Value = '  +4.69%  '
$ws.Range("D14").Value = '''25.80'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.23%  '
$ws.Range("D15").Value = '''0.0000175'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +10.10%  '
$ws.Range("D16").Value = '2.813.44'
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("D17").Value = '61.612.80'
$ws.Range("E17").Value = '  +2.37%  '
$ws.Range("D18").Value = '2.384.99'
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("D19").Value = '''11.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.11%  '
$ws.Range("D20").Value = '''4.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.85%  '
$ws.Range("D21").Value = '''323.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.52%  '
$ws.Range("D22").Value = '''6.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.23%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").Value = '  -3.80%  '
$ws.Range("D25").Value = '''64.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.53%  '
$ws.Range("D26").Value = '''9.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.00%  '
$ws.Range("D27").Value = '''542.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +10.52%  '
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").Value = '2.501.05'
$ws.Range("E29").Value = '  +1.61%  '
$ws.Range("D30").Value = '''8.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.46%  '
$ws.Range("D31").Value = '0.0₃0926'
$ws.Range("E31").Value = '  +6.00%  '
$ws.Range("E32").Value = '  +3.97%  '
$ws.Range("E33").Value = '  +3.92%  '
$ws.Range("E34").Value = '  +4.63%  '
$ws.Range("E35").Value = '  +2.20%  '
$ws.Range("D36").Value = '''5.74'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.84%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").Value = '''4.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.38%  '
$ws.Range("D39").Value = '''1.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.10%  '
$ws.Range("D40").Value = '''0.381'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.11%  '
$ws.Range("D41").Value = '''18.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.36%  '
$ws.Range("D42").Value = '''147.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.04%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '''41.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '''2.27'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.77%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''148.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.13%  '
$ws.Range("D48").Value = '''0.0531'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.04%  '
$ws.Range("D49").Value = '''20.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.91%  '
$ws.Range("D50").Value = '''0.585'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("D51").Value = '''0.0907'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.83%  '
